$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 5
$ws.Range("O6").Value = 5
